$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.223.29"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.38%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.264.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.00%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.30%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'497.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.46%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'128.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +1.17%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.54%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  -0.51%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.0952"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +0.78%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  +0.75%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +3.32%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'4.77"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +3.37%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'2.664.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.49%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'22.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +5.23%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'54.187.60"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +0.32%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'  +0.60%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'2.268.34"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +0.21%  "
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'  +2.49%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  +1.31%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'302.71"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +0.49%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'  -2.52%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  +0.54%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'60.95"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -2.15%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +0.53%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  -0.67%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'7.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +3.39%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'170.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +1.29%  "
$ws.Range("E27").ClearFormats()
$ws.Range("B28").Value = "'PEPE"
$ws.Range("B28").ClearFormats()
$ws.Range("C28").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C28").ClearFormats()
$ws.Range("D28").Value = "'0.0₃0692"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +0.61%  "
$ws.Range("E28").ClearFormats()
$ws.Range("B29").Value = "'PancakeSwap"
$ws.Range("B29").ClearFormats()
$ws.Range("C29").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C29").ClearFormats()
$ws.Range("D29").Value = "'1.60"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +0.56%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E31").Value = "'  +1.24%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +0.23%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'17.74"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +1.09%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'0.996"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +0.58%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'0.946"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +9.76%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  +0.08%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  +0.67%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  -0.19%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'1.40"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -0.17%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  +0.68%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'125.55"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -2.07%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'4.80"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -2.47%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  +2.43%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  +0.28%  "
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = "'  +0.54%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'241.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +1.62%  "
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'  -0.01%  "
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "'  +1.03%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'10.76"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +0.48%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'16.16"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -0.29%  "
$ws.Range("E50").ClearFormats()
